# Update PCB Schem + Objectif dev durable
#
# Slide 1 ("PCB Schem"):
#   - TextBox 2: "17395xx36" stays the same text (formatting refresh only).
#   - TextBox 24: "Tiny 8X RC Drone" -> split into "Tiny" + " 8X RC Drone"
#     (two runs, same visual formatting).
#
# Slide 2 ("Objectif dev durable"):
#   - Rectangle 141: resize/reposition, "Regulateur" -> "Regulateur a decoupage",
#     "LDO" -> "17395xx36" (new run picks up the default/inherited formatting).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# TextBox 2 -> "17395xx36" (text unchanged; nudge formatting so PowerPoint
# re-serialises the run/endParaRPr).
$tb2 = $s1.Shapes.Item("TextBox 2")
$tb2.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "17395xx36"

# TextBox 24 -> split "Tiny 8X RC Drone" into two runs: "Tiny" and " 8X RC Drone"
$tb24 = $s1.Shapes.Item("TextBox 24")
$tr24 = $tb24.TextFrame.TextRange
$para1 = $tr24.Paragraphs(1)
$firstPart = $para1.Characters(1, 4)
# Touching the sub-range's font forces PowerPoint to split it into its own run
# while keeping the original look (font/size/colour unchanged).
$firstPart.Font.Name = $firstPart.Font.Name
$firstPart.Font.Size = $firstPart.Font.Size
$firstPart.Font.Bold = $firstPart.Font.Bold
$firstPart.Font.Color.RGB = $firstPart.Font.Color.RGB

# ---------------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$rect = $s2.Shapes.Item("Rectangle 141")

# Reposition / resize the rectangle (EMU -> points, 1 pt = 12700 EMU).
$rect.Left = 1001159 / 12700
$rect.Width = 2625300 / 12700
$rect.Top = 1546560 / 12700
$rect.Height = 599760 / 12700

$rectTr = $rect.TextFrame.TextRange

# Paragraph 1: "Regulateur" -> "Regulateur a decoupage"
$rectTr.Paragraphs(1).Runs(1).Text = "Régulateur à découpage"

# Paragraph 2: "LDO" -> "17395xx36"
$rectTr.Paragraphs(2).Runs(1).Text = "17395xx36"
